# Update the two-digit x two-digit multiplication answers in the table
# to match the newly generated problem set (output generated at c8c62b6).
$d = $word.ActiveDocument

$d.Content.Find.Execute("38×27=1026", $true, $false, $false, $false, $false, $true, 1, $false, "90×21=1890", 2) | Out-Null
$d.Content.Find.Execute("70×84=5880", $true, $false, $false, $false, $false, $true, 1, $false, "27×14=378", 2) | Out-Null
$d.Content.Find.Execute("31×20=620", $true, $false, $false, $false, $false, $true, 1, $false, "21×99=2079", 2) | Out-Null
$d.Content.Find.Execute("91×36=3276", $true, $false, $false, $false, $false, $true, 1, $false, "63×37=2331", 2) | Out-Null
$d.Content.Find.Execute("60×63=3780", $true, $false, $false, $false, $false, $true, 1, $false, "57×88=5016", 2) | Out-Null
$d.Content.Find.Execute("15×66=990", $true, $false, $false, $false, $false, $true, 1, $false, "76×47=3572", 2) | Out-Null
$d.Content.Find.Execute("47×48=2256", $true, $false, $false, $false, $false, $true, 1, $false, "47×11=517", 2) | Out-Null
$d.Content.Find.Execute("94×76=7144", $true, $false, $false, $false, $false, $true, 1, $false, "69×95=6555", 2) | Out-Null
$d.Content.Find.Execute("24×82=1968", $true, $false, $false, $false, $false, $true, 1, $false, "40×50=2000", 2) | Out-Null
$d.Content.Find.Execute("61×49=2989", $true, $false, $false, $false, $false, $true, 1, $false, "95×33=3135", 2) | Out-Null
$d.Content.Find.Execute("66×86=5676", $true, $false, $false, $false, $false, $true, 1, $false, "12×64=768", 2) | Out-Null
$d.Content.Find.Execute("97×68=6596", $true, $false, $false, $false, $false, $true, 1, $false, "80×65=5200", 2) | Out-Null
$d.Content.Find.Execute("13×13=169", $true, $false, $false, $false, $false, $true, 1, $false, "50×15=750", 2) | Out-Null
$d.Content.Find.Execute("62×27=1674", $true, $false, $false, $false, $false, $true, 1, $false, "64×65=4160", 2) | Out-Null
$d.Content.Find.Execute("97×64=6208", $true, $false, $false, $false, $false, $true, 1, $false, "74×59=4366", 2) | Out-Null
$d.Content.Find.Execute("21×87=1827", $true, $false, $false, $false, $false, $true, 1, $false, "20×47=940", 2) | Out-Null
$d.Content.Find.Execute("14×64=896", $true, $false, $false, $false, $false, $true, 1, $false, "65×46=2990", 2) | Out-Null
$d.Content.Find.Execute("32×60=1920", $true, $false, $false, $false, $false, $true, 1, $false, "50×31=1550", 2) | Out-Null
$d.Content.Find.Execute("66×75=4950", $true, $false, $false, $false, $false, $true, 1, $false, "60×73=4380", 2) | Out-Null
$d.Content.Find.Execute("98×81=7938", $true, $false, $false, $false, $false, $true, 1, $false, "30×67=2010", 2) | Out-Null
$d.Content.Find.Execute("28×25=700", $true, $false, $false, $false, $false, $true, 1, $false, "98×81=7938", 2) | Out-Null
$d.Content.Find.Execute("77×79=6083", $true, $false, $false, $false, $false, $true, 1, $false, "57×69=3933", 2) | Out-Null
$d.Content.Find.Execute("53×31=1643", $true, $false, $false, $false, $false, $true, 1, $false, "69×35=2415", 2) | Out-Null
$d.Content.Find.Execute("72×60=4320", $true, $false, $false, $false, $false, $true, 1, $false, "20×24=480", 2) | Out-Null
$d.Content.Find.Execute("12×95=1140", $true, $false, $false, $false, $false, $true, 1, $false, "46×94=4324", 2) | Out-Null
